$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80

# Column A holds a date-looking string that must stay plain text (like the
# other rows), not get auto-converted into a date serial number.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "02/12/2026"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = 9390.84
$ws.Cells.Item($row, 3).Value = 0.2372190378971422
$ws.Cells.Item($row, 4).Value = 0.7627809621028578
$ws.Cells.Item($row, 5).Value = -312.05
$ws.Cells.Item($row, 6).Value = -37.37
$ws.Cells.Item($row, 7).Value = -23728.61
$ws.Cells.Item($row, 8).Value = -76.81
$ws.Cells.Item($row, 9).Value = -1126.02
$ws.Cells.Item($row, 10).Value = -33.58
$ws.Cells.Item($row, 11).Value = -24853.77
$ws.Cells.Item($row, 12).Value = -72.58
